$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 278, pushing the existing rows 278-342 down to 279-343
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row 278 with the new week's data
$ws.Range("A278").Value2 = 9
$ws.Range("B278").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C278").Value2 = "Metropolitana"
$ws.Range("D278").Value2 = 44889
$ws.Range("E278").Value2 = 13
$ws.Range("F278").Value2 = 300000001
$ws.Range("G278").Value2 = "Rabanito"
$ws.Range("H278").Value2 = "Sin especificar"
$ws.Range("I278").Value2 = "Primera"
$ws.Range("J278").Value2 = 11000
$ws.Range("K278").Value2 = 3000
$ws.Range("L278").Value2 = 4000
$ws.Range("M278").Value2 = 3455
$ws.Range("N278").Value2 = "`$/cien unidades (volumen en unidades)"
$ws.Range("O278").Value2 = "Provincia de Chacabuco"
$ws.Range("P278").Value2 = 35
$ws.Range("Q278").Value2 = 100
$ws.Range("R278").Value2 = "Hortaliza"
